$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty CNV (PAAD) values for rows 11 and 12
$ws.Range("C11").Value = 220
$ws.Range("C12").Value = 37

# The shared formula in C13 (=C11+C12) will recalc automatically (220+37=257)

# Update the active cell / selection to reflect the edited cell
$ws.Range("C12").Select()
